# Apply the updated crypto price/volume snapshot (and the BitcoinCash/ShibaInu row
# swap) described by the commit's XML diff. Each data row is "<CellRef><TAB><NewValue>";
# values are applied via the Excel Range object exactly like a user typing into the grid.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @'
D2	29.215.69
E2	  +0.07%  
D3	1.860.84
E3	  +0.52%  
D4	0.9995
E4	  -0.15%  
D5	0.6984
E5	  +0.55%  
D6	236.79
E6	  -0.68%  
D7	0.9987
E7	  -0.22%  
D8	0.07650
E8	  +0.98%  
D9	0.3043
E9	  -0.74%  
D10	23.26
E10	  -0.95%  
D11	0.08131
E11	  +0.37%  
D12	1.859.74
E12	  +0.32%  
D13	0.7158
E13	  -0.98%  
D14	5.139
E14	  -0.81%  
D15	89.47
E15	  +0.51%  
D16	29.227.28
E16	  +0.03%  
D17	5.742
E17	  -0.79%  
D18	13.15
E18	  +0.48%  
B19	ShibaInu
C19	https://coinranking.com/coin/xz24e0BjL+shibainu-shib
D19	0.000007702
E19	  -0.19%  
B20	BitcoinCash
C20	https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch
D20	237.21
E20	  -1.53%  
D21	0.9998
E21	  +0.00%  
D22	2.112.22
E22	  +0.80%  
D23	0.9994
E23	  -0.19%  
D24	7.484
E24	  -1.60%  
D25	9.007
E25	  -0.09%  
D26	161.71
E26	  +0.10%  
D27	0.1454
E27	  +0.01%  
D28	18.04
E28	  -0.09%  
D29	1.970
E29	  +1.95%  
E30	  +0.84%  
D31	4.443
E31	  +0.42%  
D32	1.478
E32	  -1.57%  
D33	3.992
E33	  -1.18%  
D34	0.05181
E34	  -0.78%  
D35	1.166
E35	  -1.96%  
D36	0.7067
E36	  -0.17%  
D37	0.9987
E37	  -0.03%  
D38	2.650
E38	  -0.51%  
D39	0.01853
E39	  -0.37%  
D40	2.717
E40	  +0.89%  
D41	0.9318
E41	  -0.34%  
D42	1.135.58
E42	  +8.53%  
E43	  -0.48%  
D44	70.72
E44	  +1.76%  
D45	5.872
E45	  -1.29%  
D46	0.9979
E46	  -0.31%  
D47	103.43
E47	  +1.06%  
E48	  +3.69%  
D49	2.010.91
E49	  +0.80%  
D50	9.165
E50	  -1.25%  
D51	6.948
E51	  -3.95%  
'@

$rows = $updates -split "`r?`n" | Where-Object { $_.Length -gt 0 }

foreach ($row in $rows) {
    $parts = $row -split "`t", 2
    $cellRef = $parts[0]
    $newValue = $parts[1]

    $range = $ws.Range($cellRef)

    # These columns hold free-text price/volume strings (e.g. "29.215.69", "0.9995",
    # "  +0.07%  "). Excel auto-detects plain numeric-looking text as a number on
    # assignment, which would corrupt values like "0.9995" (-> 0.99950000000000006)
    # or lose the original "29.247.01"-style grouping. Force the cell to Text format
    # first whenever the new value would otherwise parse as a plain number, so it's
    # stored verbatim as a string, matching the source data.
    if ($newValue -match '^-?\d+(\.\d+)?$') {
        $range.NumberFormat = "@"
    }

    $range.Value = $newValue
}
